$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; E=2; G=1.7502375; H=3.500475; I=0.06096295440307654; J=0.05231488010959406; K=2; M=2.019046; N=4.038092; O=0.003975353327590414; P=0.002707057536954368; Q=3.533810023425; R=14.1352400937; S=0.000242349283646013; T=0.0001416193904955408 },
    @{ Row=3; E=2; G=1.7502375; H=3.500475; I=0.06096295440307654; J=0.05231488010959406; K=3; M=139.6948166666666; N=419.0844499999999; O=0.2750488370661026; P=0.2809459811695414; Q=244.4991066856249; R=1466.99464011375; S=0.01676778971268004; T=0.01469765532215683 },
    @{ Row=4; E=2; G=1.7502375; H=3.500475; I=0.06096295440307654; J=0.05231488010959406; K=3; M=186.3548536666667; N=559.064561; O=0.3669190239530987; P=0.3747858972750337; Q=326.1652531944125; R=1956.991519166475; S=0.02236846772687411; T=0.01960687928271002 },
    @{ Row=5; E=2; G=1.7502375; H=3.500475; I=0.06096295440307654; J=0.05231488010959406; K=3; M=143.6051993333333; N=430.815598; O=0.2827480934208787; P=0.2888103122968479; Q=251.343205068175; R=1508.05923040905; S=0.01723715912677385; T=0.01510907686222402 },
    @{ Row=6; E=2; G=1.7502375; H=3.500475; I=0.06096295440307654; J=0.05231488010959406; K=3; M=6.253715333333335; N=18.761146; O=0.01231310631861279; P=0.0125771036621259; Q=10.945487090725; R=65.67292254435002; S=0.000750643339061825; T=0.000657969670210053 },
    @{ Row=7; E=2; G=1.7502375; H=3.500475; I=0.06096295440307654; J=0.05231488010959406; K=2; M=29.9633245; N=59.926649; O=0.05899558591371687; P=0.04017364805949665; Q=52.44293416456875; R=209.771736658275; S=0.003596545214040706; T=0.002101679581797593 },
    @{ Row=8; E=3; G=3.198312; H=9.594936; I=0.1114011947651748; J=0.1433970893947902; K=2; M=2.019046; N=4.038092; O=0.003975353327590414; P=0.002707057536954368; Q=6.457539050352; R=38.745234302112; S=0.0004428591103072853; T=0.000388184171623486 },
    @{ Row=9; E=3; G=3.198312; H=9.594936; I=0.1114011947651748; J=0.1433970893947902; K=3; M=139.6948166666666; N=419.0844499999999; O=0.2750488370661026; P=0.2809459811695414; Q=446.7876084827999; R=4021.0884763452; S=0.03064076906793571; T=0.04028683597687577 },
    @{ Row=10; E=3; G=3.198312; H=9.594936; I=0.1114011947651748; J=0.1433970893947902; K=3; M=186.3548536666667; N=559.064561; O=0.3669190239530987; P=0.3747858972750337; Q=596.0209647403441; R=5364.188682663096; S=0.04087521765044697; T=0.05374320681545464 },
    @{ Row=11; E=3; G=3.198312; H=9.594936; I=0.1114011947651748; J=0.1433970893947902; K=3; M=143.6051993333333; N=430.815598; O=0.2827480934208787; P=0.2888103122968479; Q=459.294232290192; R=4133.648090611729; S=0.03149847542466113; T=0.04141455817056837 },
    @{ Row=12; E=3; G=3.198312; H=9.594936; I=0.1114011947651748; J=0.1433970893947902; K=3; M=6.253715333333335; N=18.761146; O=0.01231310631861279; P=0.0125771036621259; Q=20.00133279518401; R=180.011995156656; S=0.001371694755164087; T=0.001803520058165411 },
    @{ Row=13; E=3; G=3.198312; H=9.594936; I=0.1114011947651748; J=0.1433970893947902; K=2; M=29.9633245; N=59.926649; O=0.05899558591371687; P=0.04017364805949665; Q=95.832060308244; R=574.992361849464; S=0.006572178756659573; T=0.005760784202102479 },
    @{ Row=14; E=3; G=2.564619666666667; H=7.693859; I=0.0893289006778985; J=0.1149853408937705; K=2; M=2.019046; N=4.038092; O=0.003975353327590414; P=0.002707057536954368; Q=5.178085079504666; R=31.068510477028; S=0.0003551139425598774; T=0.0003112719337057487 },
    @{ Row=15; E=3; G=2.564619666666667; H=7.693859; I=0.0893289006778985; J=0.1149853408937705; K=3; M=139.6948166666666; N=419.0844499999999; O=0.2750488370661026; P=0.2809459811695414; Q=358.2640741547277; R=3224.376667392549; S=0.02456981024784936; T=0.03230466941751455 },
    @{ Row=16; E=3; G=2.564619666666667; H=7.693859; I=0.0893289006778985; J=0.1149853408937705; K=3; M=186.3548536666667; N=559.064561; O=0.3669190239530987; P=0.3747858972750337; Q=477.9293226923222; R=4301.363904230899; S=0.03277647304753781; T=0.0430948841603474 },
    @{ Row=17; E=3; G=2.564619666666667; H=7.693859; I=0.0893289006778985; J=0.1149853408937705; K=3; M=143.6051993333333; N=430.815598; O=0.2827480934208787; P=0.2888103122968479; Q=368.2927184458536; R=3314.634466012682; S=0.02525757635405884; T=0.03320895221308938 },
    @{ Row=18; E=3; G=2.564619666666667; H=7.693859; I=0.0893289006778985; J=0.1149853408937705; K=3; M=6.253715333333335; N=18.761146; O=0.01231310631861279; P=0.0125771036621259; Q=16.03840133360156; R=144.345612002414; S=0.001099916251371766; T=0.001446182552045836 },
    @{ Row=19; E=3; G=2.564619666666667; H=7.693859; I=0.0893289006778985; J=0.1149853408937705; K=2; M=29.9633245; N=59.926649; O=0.05899558591371687; P=0.04017364805949665; Q=76.84453129141517; R=461.067187748491; S=0.005270010834520843; T=0.004619380617067584 },
    @{ Row=20; E=3; G=2.636255; H=7.908765000000001; I=0.09182404865618671; J=0.1181971283297135; K=2; M=2.019046; N=4.038092; O=0.003975353327590414; P=0.002707057536954368; Q=5.32272011273; R=31.93632067638; S=0.0003650330373781959; T=0.0003199664270913135 },
    @{ Row=21; E=3; G=2.636255; H=7.908765000000001; I=0.09182404865618671; J=0.1181971283297135; K=3; M=139.6948166666666; N=419.0844499999999; O=0.2750488370661026; P=0.2809459811695414; Q=368.2711589115833; R=3314.44043020425; S=0.02525609779758537; T=0.03320700819001355 },
    @{ Row=22; E=3; G=2.636255; H=7.908765000000001; I=0.09182404865618671; J=0.1181971283297135; K=3; M=186.3548536666667; N=559.064561; O=0.3669190239530987; P=0.3747858972750337; Q=491.2789147530184; R=4421.510232777166; S=0.03369199030834987; T=0.04429861679638396 },
    @{ Row=23; E=3; G=2.636255; H=7.908765000000001; I=0.09182404865618671; J=0.1181971283297135; K=3; M=143.6051993333333; N=430.815598; O=0.2827480934208787; P=0.2888103122968479; Q=378.5799247684967; R=3407.21932291647; S=0.02596307468772279; T=0.03413654954549516 },
    @{ Row=24; E=3; G=2.636255; H=7.908765000000001; I=0.09182404865618671; J=0.1181971283297135; K=3; M=6.253715333333335; N=18.761146; O=0.01231310631861279; P=0.0125771036621259; Q=16.48638831607667; R=148.37749484469; S=0.0011306392737091; T=0.001486577535568405 },
    @{ Row=25; E=3; G=2.636255; H=7.908765000000001; I=0.09182404865618671; J=0.1181971283297135; K=2; M=29.9633245; N=59.926649; O=0.05899558591371687; P=0.04017364805949665; Q=78.9909640297475; R=473.945784178485; S=0.005417213551441381; T=0.00474840983516107 },
    @{ Row=26; E=3; G=1.092757; H=3.278271; I=0.03806208881060012; J=0.04899402347731637; K=2; M=2.019046; N=4.038092; O=0.003975353327590414; P=0.002707057536954368; Q=2.206326649822; R=13.237959898932; S=0.000151310251408261; T=0.0001326296405199886 },
    @{ Row=27; E=3; G=1.092757; H=3.278271; I=0.03806208881060012; J=0.04899402347731637; K=3; M=139.6948166666666; N=419.0844499999999; O=0.2750488370661026; P=0.2809459811695414; Q=152.6524887762166; R=1373.87239898595; S=0.01046893326366228; T=0.0137646739972782 },
    @{ Row=28; E=3; G=1.092757; H=3.278271; I=0.03806208881060012; J=0.04899402347731637; K=3; M=186.3548536666667; N=559.064561; O=0.3669190239530987; P=0.3747858972750337; Q=203.6405708282257; R=1832.765137454031; S=0.01396570447600155; T=0.01836226905006008 },
    @{ Row=29; E=3; G=1.092757; H=3.278271; I=0.03806208881060012; J=0.04899402347731637; K=3; M=143.6051993333333; N=430.815598; O=0.2827480934208787; P=0.2888103122968479; Q=156.9255868078953; R=1412.330281271058; S=0.01076198304281334; T=0.01414997922116284 },
    @{ Row=30; E=3; G=1.092757; H=3.278271; I=0.03806208881060012; J=0.04899402347731637; K=3; M=6.253715333333335; N=18.761146; O=0.01231310631861279; P=0.0125771036621259; Q=6.833791206507335; R=61.504120858566; S=0.0004686625462334013; T=0.0006162029120988382 },
    @{ Row=31; E=3; G=1.092757; H=3.278271; I=0.03806208881060012; J=0.04899402347731637; K=2; M=29.9633245; N=59.926649; O=0.05899558591371687; P=0.04017364805949665; Q=32.7426325906465; R=196.455795543879; S=0.00224549523048128; T=0.001968268656196424 },
    @{ Row=32; E=2; G=17.4676725; H=34.935345; I=0.6084208126870634; J=0.5221115377948153; K=2; M=2.019046; N=4.038092; O=0.003975353327590414; P=0.002707057536954368; Q=35.268034290435; R=141.07213716174; S=0.002418687702290782; T=0.00141338597351829 },
    @{ Row=33; E=2; G=17.4676725; H=34.935345; I=0.6084208126870634; J=0.5221115377948153; K=3; M=139.6948166666666; N=419.0844499999999; O=0.2750488370661026; P=0.2809459811695414; Q=2440.143307480874; R=14640.85984488525; S=0.1673454369763898; T=0.1466851382657025 },
    @{ Row=34; E=2; G=17.4676725; H=34.935345; I=0.6084208126870634; J=0.5221115377948153; K=3; M=186.3548536666667; N=559.064561; O=0.3669190239530987; P=0.3747858972750337; Q=3255.185552634758; R=19531.11331580854; S=0.2232411707438884; T=0.1956800411700775 },
    @{ Row=35; E=2; G=17.4676725; H=34.935345; I=0.6084208126870634; J=0.5221115377948153; K=3; M=143.6051993333333; N=430.815598; O=0.2827480934208787; P=0.2888103122968479; Q=2508.448591251885; R=15050.69154751131; S=0.1720298247848487; T=0.1507911962843081 },
    @{ Row=36; E=2; G=17.4676725; H=34.935345; I=0.6084208126870634; J=0.5221115377948153; K=3; M=6.253715333333335; N=18.761146; O=0.01231310631861279; P=0.0125771036621259; Q=109.237851350895; R=655.42710810537; S=0.007491550153072606; T=0.006566650934037359 },
    @{ Row=37; E=2; G=17.4676725; H=34.935345; I=0.6084208126870634; J=0.5221115377948153; K=2; M=29.9633245; N=59.926649; O=0.05899558591371687; P=0.04017364805949665; Q=523.3895393772262; R=2093.558157508905; S=0.03589414232657308; T=0.0209751251671715 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
    $ws.Cells.Item($r.Row, 10).Value = $r.J
    $ws.Cells.Item($r.Row, 11).Value = $r.K
    $ws.Cells.Item($r.Row, 13).Value = $r.M
    $ws.Cells.Item($r.Row, 14).Value = $r.N
    $ws.Cells.Item($r.Row, 15).Value = $r.O
    $ws.Cells.Item($r.Row, 16).Value = $r.P
    $ws.Cells.Item($r.Row, 17).Value = $r.Q
    $ws.Cells.Item($r.Row, 18).Value = $r.R
    $ws.Cells.Item($r.Row, 19).Value = $r.S
    $ws.Cells.Item($r.Row, 20).Value = $r.T
}
